$d = $word.ActiveDocument

# 1. Remove the motion-related sentences that were appended after the
#    arraignment sentence (4 runs: the motion text, the court's ruling,
#    and two trailing single-space runs). The empty run immediately
#    preceding them is left untouched.
$rng = $d.Content
$found = $rng.Find.Execute("Counsel for the State of Ohio made a motion", $true, $false, $false,
                            $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $para = $rng.Paragraphs(1)
    $paraEnd = $para.Range.End
    $delRange = $d.Range($rng.Start, $paraEnd - 1)
    $delRange.Delete()
}

# 2. Update the offense/charge table cell text ("DUS UCM - AMENDED" -> "DUS UCM").
$d.Content.Find.Execute("DUS UCM - AMENDED", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DUS UCM", 2) | Out-Null

# 3. Update the four plea cells from "No Contest" to "Guilty".
$d.Content.Find.Execute("No Contest", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Guilty", 2) | Out-Null
